# Refresh the cryptocurrency "Price" (D) and "Volume(1h)" (E) columns for
# rows 2-51, matching the data snapshot from the
# "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.936.48"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "3.653.31"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'597.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.29%  "
$ws.Range("D6").Value = "'189.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.20%  "
$ws.Range("D7").Value = "'0.620"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").Value = "'0.701"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("E10").Value = "  +6.22%  "
$ws.Range("E11").Value = "  -6.34%  "
$ws.Range("E12").Value = "  -6.83%  "
$ws.Range("D13").Value = "'10.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("D14").Value = "4.237.18"
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("D15").Value = "3.658.54"
$ws.Range("E15").Value = "  -2.00%  "
$ws.Range("D16").Value = "'0.126"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "'18.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.55%  "
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").Value = "67.731.28"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("E20").Value = "  -3.19%  "
$ws.Range("D21").Value = "'400.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("E22").Value = "  -2.24%  "
$ws.Range("D23").Value = "'87.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("D24").Value = "'11.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("E25").Value = "  -2.35%  "
$ws.Range("E26").Value = "  -2.40%  "
$ws.Range("D27").Value = "'6.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").Value = "'3.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.39%  "
$ws.Range("D29").Value = "'9.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.52%  "
$ws.Range("D30").Value = "'31.78"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("D31").Value = "'7.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.23%  "
$ws.Range("D32").Value = "'12.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.83%  "
$ws.Range("D33").Value = "'44.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.79%  "
$ws.Range("D34").Value = "'65.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("D36").Value = "'606.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("E40").Value = "  -14.22%  "
$ws.Range("D41").Value = "'0.135"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("D42").Value = "'2.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.07%  "
$ws.Range("D43").Value = "'0.0426"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("D44").Value = "'2.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.53%  "
$ws.Range("D45").Value = "'0.136"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("D46").Value = "2.774.71"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").Value = "'3.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("D48").Value = "'143.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.62%  "
$ws.Range("E49").Value = "  -5.64%  "
$ws.Range("E50").Value = "  -3.47%  "
$ws.Range("D51").Value = "'2.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -15.02%  "
